$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "73.211.63"
$ws.Range("D3").Value = "3.993.02"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'611.23"
$ws.Range("D6").Value = "'166.48"
$ws.Range("E6").Value = "  +11.36%  "
$ws.Range("D7").Value = "'0.684"
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").Value = "'56.89"
$ws.Range("E11").Value = "  +6.43%  "
$ws.Range("D12").Value = "'0.0000337"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "'11.12"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "4.627.79"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").Value = "4.003.87"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("D17").Value = "'14.21"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "'20.58"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "73.043.40"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").Value = "'439.02"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'4.94"
$ws.Range("E22").Value = "  +15.59%  "
$ws.Range("D23").Value = "'95.99"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").Value = "'3.37"
$ws.Range("E24").Value = "  -3.95%  "
$ws.Range("D25").Value = "'14.23"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").Value = "'4.09"
$ws.Range("E26").Value = "  -6.33%  "
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'10.54"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "'5.97"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "'36.09"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").Value = "'7.65"
$ws.Range("E31").Value = "  -8.69%  "
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "'0.0000103"
$ws.Range("E34").Value = "  +18.29%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'72.14"
$ws.Range("E35").Value = "  +8.15%  "
$ws.Range("D36").Value = "'48.15"
$ws.Range("E36").Value = "  -4.26%  "
$ws.Range("D37").Value = "'634.10"
$ws.Range("E37").Value = "  -6.62%  "
$ws.Range("D38").Value = "'0.432"
$ws.Range("E38").Value = "  -6.20%  "
$ws.Range("D39").Value = "'3.42"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("D42").Value = "'11.09"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'3.26"
$ws.Range("E44").Value = "  -5.44%  "
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("D46").Value = "'0.149"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("E47").Value = "  +4.33%  "
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("E49").Value = "  +29.68%  "
$ws.Range("D50").Value = "2.873.33"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000282"
$ws.Range("E51").Value = "  +1.97%  "
